$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Record a row selection on the existing DXX sheet (cosmetic view
#    state captured in the diff: selection moved from K11 to A6:XFD7,
#    i.e. the user selected rows 6:7 before adding the new sheet).
# ------------------------------------------------------------------
$dxx = $wb.Worksheets.Item("DXX")
$dxx.Activate()
$dxx.Rows("6:7").Select()

# ------------------------------------------------------------------
# 2) Insert a brand new worksheet "DXXAG" right after "DXX" (this is
#    what pushes FASTQX..TRIGLY/sheet12..29 down by one slot).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $dxx)
$ws.Name = "DXXAG"

# ------------------------------------------------------------------
# 3) Header row (row 1) -- same year columns used across every other
#    variable sheet in the workbook.
# ------------------------------------------------------------------
$headers = @("variable","1999-2000","2001-2002","2003-2004","2005-2006","2007-2008","2009-2010","2011-2012","2013-2014","2015-2016","2017-2018","2017-Mar2020","2021-2023")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}
$ws.Range("A1:M1").Font.Bold = $true

# ------------------------------------------------------------------
# 4) Row 2 -- respondent id / SEQN, identical to every other sheet.
# ------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "respondentid"
for ($col = 2; $col -le 13; $col++) {
    $ws.Cells.Item(2, $col).Value = "SEQN"
}

# ------------------------------------------------------------------
# 5) Row labels first -- visceral_fat (row 4) was entered before
#    subcutaneous_fat (row 3), matching the shared-string insertion
#    order of the target workbook.
# ------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "visceral_fat"
$ws.Cells.Item(3, 1).Value = "subcutaneous_fat"

# ------------------------------------------------------------------
# 6) Row 4 -- visceral_fat / DXXVFATM (2011-2012 .. 2017-2018,
#    columns H:K), entered before row 3's values.
# ------------------------------------------------------------------
$ws.Cells.Item(4, 8).Value = "DXXVFATM"
$ws.Cells.Item(4, 9).Value = "DXXVFATM"
$ws.Cells.Item(4, 10).Value = "DXXVFATM"
$ws.Cells.Item(4, 11).Value = "DXXVFATM"

# ------------------------------------------------------------------
# 7) Row 3 -- subcutaneous_fat / DXXSATM (2011-2012 .. 2017-2018,
#    columns H:K).
# ------------------------------------------------------------------
$ws.Cells.Item(3, 8).Value = "DXXSATM"
$ws.Cells.Item(3, 9).Value = "DXXSATM"
$ws.Cells.Item(3, 10).Value = "DXXSATM"
$ws.Cells.Item(3, 11).Value = "DXXSATM"

# ------------------------------------------------------------------
# 8) Leave the new sheet active with the same cell selection found in
#    the target workbook.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("G6").Select()
